# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column D = Price, Column E = Volume(1h). Values that read as plain numbers
# are prefixed with a leading apostrophe so Excel stores them as text (matching
# the workbook's original inlineStr text cells, e.g. "1.000" must stay "1.000"
# rather than becoming the number 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.014.02'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '1.829.16'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''324.06'
$ws.Range("E5").Value = '  -3.72%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '''0.4616'
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").Value = '''0.3853'
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").Value = '''0.07825'
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("D10").Value = '''0.9579'
$ws.Range("E10").Value = '  -2.53%  '
$ws.Range("D11").Value = '''21.84'
$ws.Range("E11").Value = '  -2.23%  '
$ws.Range("D12").Value = '1.846.58'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("D13").Value = '''5.662'
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").Value = '''6.871'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").Value = '''0.06865'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '''88.19'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '''1.002'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '''0.000009897'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").Value = '''16.64'
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '28.024.34'
$ws.Range("E21").Value = '  -2.03%  '
$ws.Range("D22").Value = '''5.283'
$ws.Range("E22").Value = '  -2.32%  '
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("D24").Value = '''2.075'
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("D25").Value = '2.058.80'
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("D26").Value = '''154.39'
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("D27").Value = '''19.12'
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("D28").Value = '''5.673'
$ws.Range("E28").Value = '  -6.63%  '
$ws.Range("E29").Value = '  -3.44%  '
$ws.Range("D30").Value = '''118.27'
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").Value = '''0.9354'
$ws.Range("E31").Value = '  -3.90%  '
$ws.Range("D32").Value = '''0.09211'
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").Value = '''5.247'
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("D34").Value = '''1.315'
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("D35").Value = '''3.306'
$ws.Range("E35").Value = '  -5.14%  '
$ws.Range("D36").Value = '''0.05817'
$ws.Range("E36").Value = '  -5.59%  '
$ws.Range("D37").Value = '''0.02116'
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("D39").Value = '''7.710'
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").Value = '''0.5571'
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("D41").Value = '''9.857'
$ws.Range("E41").Value = '  -3.29%  '
$ws.Range("D42").Value = '''0.1755'
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("D43").Value = '''0.07367'
$ws.Range("E43").Value = '  +3.19%  '
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("E45").Value = '  -2.70%  '
$ws.Range("D46").Value = '''1.128'
$ws.Range("E46").Value = '  -10.05%  '
$ws.Range("D47").Value = '''2.089'
$ws.Range("E47").Value = '  -11.35%  '
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("D49").Value = '''112.95'
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").Value = '''1.000'
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '''1.021'
$ws.Range("E51").Value = '  -0.16%  '
